$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / recalculated means
$ws.Range("F2").Value = -6
$ws.Range("F6").Value = -3
$ws.Range("F9").Value = 3
$ws.Range("F13").Value = 0
$ws.Range("F17").Value = 6
$ws.Range("F20").Value = -3
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = -5
$ws.Range("F26").Value = -3
$ws.Range("F30").Value = -4
$ws.Range("F33").Value = 1
$ws.Range("F36").Value = -2
